# Fix ticker typo: "APPL" should always have been "AAPL".
# Every cell in the workbook that contained the string "APPL" is updated
# to "AAPL" (the stock ticker symbol was simply misspelled).

$wb = $excel.ActiveWorkbook

$ws_rsu = $wb.Worksheets.Item("rsu")
$ws_rsu.Range("B6").Value = "AAPL"

$ws_dividends = $wb.Worksheets.Item("dividends")
$ws_dividends.Range("B3").Value = "AAPL"

$ws_sell_orders = $wb.Worksheets.Item("sell_orders")
$ws_sell_orders.Range("B6").Value = "AAPL"
$ws_sell_orders.Range("B7").Value = "AAPL"
